# Loan RBI, Variable Instalments
#
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet. This shifts the old N/O/P columns ("Late" / "Heading" / "Outstanding")
# one column to the right (-> O/P/Q), matching the updated report layout.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment Schedule")
[void]$ws.Columns("N").Insert()

# Update the remembered selection / scroll position on this sheet.
$ws.Range("M18").Select() | Out-Null

# The "Transactions" sheet stays the active tab, but its remembered
# selection moves too.
$ws4 = $wb.Worksheets.Item("Transactions")
$ws4.Activate() | Out-Null
$ws4.Range("I8").Select() | Out-Null
